$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4797
$ws.Range("K3").Value = 4939
$ws.Range("K4").Value = 1019
$ws.Range("K5").Value = 352
$ws.Range("K6").Value = 5553
$ws.Range("K7").Value = 16660

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 146
$ws.Range("K6").Value = 125
$ws.Range("K7").Value = 490
$ws.Range("K8").Value = 1114
$ws.Range("K9").Value = 73
$ws.Range("K11").Value = 326
$ws.Range("K21").Value = 49
$ws.Range("K29").Value = 891
$ws.Range("K30").Value = 62
$ws.Range("K33").Value = 709
$ws.Range("K34").Value = 89
$ws.Range("K36").Value = 215
$ws.Range("K37").Value = 562
$ws.Range("K42").Value = 620
$ws.Range("K47").Value = 111
$ws.Range("K49").Value = 93
$ws.Range("K51").Value = 210
$ws.Range("K52").Value = 437
$ws.Range("K55").Value = 190
$ws.Range("K57").Value = 58
$ws.Range("K63").Value = 53
$ws.Range("K67").Value = 637
$ws.Range("K68").Value = 44
$ws.Range("K69").Value = 38
$ws.Range("K70").Value = 29
$ws.Range("K72").Value = 77
$ws.Range("K73").Value = 141
$ws.Range("K75").Value = 57
$ws.Range("K76").Value = 227
$ws.Range("K77").Value = 121
$ws.Range("K78").Value = 196
$ws.Range("K80").Value = 60
$ws.Range("K83").Value = 365
$ws.Range("K84").Value = 124
$ws.Range("K85").Value = 759
$ws.Range("K89").Value = 234
$ws.Range("K91").Value = 178
$ws.Range("K94").Value = 218
$ws.Range("K96").Value = 179
$ws.Range("K97").Value = 134
$ws.Range("K98").Value = 83
$ws.Range("K99").Value = 282
$ws.Range("K101").Value = 16660

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 170
$ws.Range("K3").Value = 156
$ws.Range("K7").Value = 490

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 108
$ws.Range("K7").Value = 326

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 74
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K6").Value = 179
$ws.Range("K7").Value = 759

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 119
$ws.Range("K6").Value = 163
$ws.Range("K7").Value = 437

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 312
$ws.Range("K6").Value = 374
$ws.Range("K7").Value = 1114

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 132
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 365

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 194
$ws.Range("K3").Value = 266
$ws.Range("K5").Value = 16
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 709

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 158
$ws.Range("K3").Value = 186
$ws.Range("K6").Value = 167
$ws.Range("K7").Value = 562

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 221
$ws.Range("K7").Value = 637

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 258
$ws.Range("K3").Value = 318
$ws.Range("K6").Value = 246
$ws.Range("K7").Value = 891

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 43
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 194
$ws.Range("K6").Value = 234
$ws.Range("K7").Value = 620

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 55
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 83
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 66
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 60
